$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 from 2 to 4 (Paket 1 score for row 2)
$ws.Range("B2").Value = 4

# Fill in B4 (Paket 1 score for Gian, row 4)
$ws.Range("B4").Value = 4

# Update the active selection to E5
$ws.Range("E5").Select()
